# StatementXLS: currency list import was implemented.
# Populate the "Позиция по ДС" (currency position) block on Лист1 with the
# imported currency-list figures (rate / T+N day / amount columns G:I for
# rows 13-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G13").Value = 55.85
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 120

$ws.Range("G14").Value = 0.12
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 80

$ws.Range("G15").Value = 0.12
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 80

$ws.Range("G16").Value = 55.97
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0

$ws.Range("G17").Value = 55.97
$ws.Range("H17").Value = 2
$ws.Range("I17").Value = 200

$ws.Range("G18").Value = 55.97
$ws.Range("H18").Value = 2
$ws.Range("I18").Value = 200

$ws.Range("G19").Value = 55.97
$ws.Range("H19").Value = 2
$ws.Range("I19").Value = 200

$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 2
$ws.Range("I20").Value = 200

$ws.Range("H21").Value = 70

# Reflect where the user ended up after entering the import - selection
# moves from D53 to J16 on the active (first) sheet.
$ws.Activate()
$ws.Range("J16").Select()
